# Change the table style applied to the table on slide 5 ("B1- TYPES OF
# FINANCIAL DOCUMENTS") away from the deck's custom "Table_0" style
# ({4636919D-575A-4275-BD41-DA7EBBCA8D00}) to PowerPoint's built-in
# "No Style, Table Grid" style ({2A415BDF-C217-4F0E-B14C-3F6E53E8C95E}).
$oldStyleId = "{4636919D-575A-4275-BD41-DA7EBBCA8D00}"
$newStyleId = "{2A415BDF-C217-4F0E-B14C-3F6E53E8C95E}"

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shp = $slide.Shapes.Item($shapeIdx)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
